# Append a new data row (row 48) to each of the 4 worksheets, mirroring the
# structure/style of the last existing row (row 47) and updating the values
# that differ for the new day's record.

$wb = $excel.ActiveWorkbook

$newDate = 45834.46344907407

# ---- Sheet 1: MID_LFT_#1 ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A47:I47").Copy($ws1.Range("A48:I48"))
$ws1.Cells.Item(48,1).Value = $newDate
$ws1.Cells.Item(48,4).Value = "0x01,0x6C"
$ws1.Cells.Item(48,8).Value = 364

# ---- Sheet 2: MID_LFT_#2 ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A47:I47").Copy($ws2.Range("A48:I48"))
$ws2.Cells.Item(48,1).Value = $newDate

# ---- Sheet 3: MID_PLT_#1 ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A47:I47").Copy($ws3.Range("A48:I48"))
$ws3.Cells.Item(48,1).Value = $newDate

# ---- Sheet 4: MID_PLT_#2 ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A47:I47").Copy($ws4.Range("A48:I48"))
$ws4.Cells.Item(48,1).Value = $newDate
$ws4.Cells.Item(48,4).Value = "0x00,0x7E"
$ws4.Cells.Item(48,8).Value = 126
